$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 91: correct the recorded date/time in column A (was 08:00:00, should be 07:00:00) ---
$ws.Cells.Item(91, 1).Value = 45446.2916666667

# --- Row 92: new observation appended by the R script for ticker KK.MI ---
# Give A92 the same date-time number format as the rest of column A before writing the value.
$ws.Cells.Item(91, 1).Copy()
$ws.Cells.Item(92, 1).PasteSpecial(-4122)   # xlPasteFormats

$ws.Cells.Item(92, 1).Value = 45447.4500115741   # date
$ws.Cells.Item(92, 2).Value = 3600               # volume
$ws.Cells.Item(92, 3).Value = 2                  # high
$ws.Cells.Item(92, 4).Value = 1.98000001907349   # low
$ws.Cells.Item(92, 5).Value = 2                  # open
$ws.Cells.Item(92, 6).Value = 1.98000001907349   # close

# adj_close (column G) is stored as text in this sheet, same as every other data row.
# Force text typing via a temporary "@" number format, then drop back to the default
# (unstyled) cell style so no stray "s" attribute is left on the cell.
$ws.Cells.Item(92, 7).NumberFormat = "@"
$ws.Cells.Item(92, 7).Value = "1.98000001907349"
$ws.Cells.Item(92, 7).Style = "Normal"

$ws.Cells.Item(92, 8).Value = "KK.MI"            # ticker
